# Update the "data" worksheet: remove the auto-increment "ADD_HORIZON_ID"
# column (column A), shifting Index / UsageExample / Descript left by one
# column. This matches the author's commit "update excel tables" which
# dropped the extra id column that was added for the insert_tables script.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Columns.Item(1).Delete()
